$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / label updates ---

# Header rename: "Feedstock & Reactor Conditions" -> "Feedstock & Operating Conditions"
$ws.Range("A1").Value = "Feedstock & Operating Conditions"

# Two label pairs were reordered in the source data (rows keep their numeric
# data; only the feedstock/condition label shown in column A changes):
# Row 13/14: "BSG, 250°C, 1hr" <-> "SRU-DCW-BSG, 250°C, 1hr"
$ws.Range("A13").Value = "SRU-DCW-BSG, 250°C, 1hr"
$ws.Range("A14").Value = "BSG, 250°C, 1hr"

# Row 22/23: "SRU-BSG, 190°C, 1hr" <-> "SRU-BSG, 220°C, 1hr"
$ws.Range("A22").Value = "SRU-BSG, 220°C, 1hr"
$ws.Range("A23").Value = "SRU-BSG, 190°C, 1hr"

# --- Numeric data updates (Mean / StdDev / Coefficient of Variation) ---

$ws.Range("B2").Value = 0.005096526523521016
$ws.Range("C2").Value = 0.0000953426697483441
$ws.Range("D2").Value = 1.870738223539649
$ws.Range("B3").Value = 0.004936975695214136
$ws.Range("C3").Value = 0.00008980447916251954
$ws.Range("D3").Value = 1.819018052885601
$ws.Range("B4").Value = 0.004270416550282341
$ws.Range("C4").Value = 0.00007734084601477531
$ws.Range("D4").Value = 1.811084354514828
$ws.Range("B5").Value = 0.00531261640172687
$ws.Range("C5").Value = 0.00008336004097382203
$ws.Range("D5").Value = 1.569095802714568
$ws.Range("B6").Value = 0.005603024929890132
$ws.Range("C6").Value = 0.00008584972966013805
$ws.Range("D6").Value = 1.532203242612048
$ws.Range("B7").Value = 0.004197613893729362
$ws.Range("C7").Value = 0.00006399399532739826
$ws.Range("D7").Value = 1.524532673741055
$ws.Range("B8").Value = 0.001566309759733813
$ws.Range("C8").Value = 0.00002099603016950449
$ws.Range("D8").Value = 1.340477516597525
$ws.Range("B9").Value = 0.001385352065613095
$ws.Range("C9").Value = 0.00001844779331362705
$ws.Range("D9").Value = 1.331632136807251
$ws.Range("B10").Value = 0.001149397800946305
$ws.Range("C10").Value = 0.00001476354549889915
$ws.Range("D10").Value = 1.284459173903435
$ws.Range("B11").Value = 0.0009384660007817819
$ws.Range("C11").Value = 0.000009324978256749721
$ws.Range("D11").Value = 0.9936404993874706
$ws.Range("B12").Value = 0.001047840078470718
$ws.Range("C12").Value = 0.00001041156395046639
$ws.Range("D12").Value = 0.9936214661364799
$ws.Range("B13").Value = 0.001038530027481945
$ws.Range("C13").Value = 0.00001022930702233199
$ws.Range("D13").Value = 0.9849794181815152
$ws.Range("B14").Value = 0.001035792639066947
$ws.Range("C14").Value = 0.00001019136828711756
$ws.Range("D14").Value = 0.983919744428581
$ws.Range("B15").Value = 0.0009253612996423156
$ws.Range("C15").Value = 0.000009065239146531024
$ws.Range("D15").Value = 0.9796432107151072
$ws.Range("B16").Value = 0.0008041702477353734
$ws.Range("C16").Value = 0.000007716474712294885
$ws.Range("D16").Value = 0.9595573492087373
$ws.Range("B17").Value = 0.0007946913096355002
$ws.Range("C17").Value = 0.00000754872474955552
$ws.Range("D17").Value = 0.9498939598342759
$ws.Range("B18").Value = 0.0007886472579763691
$ws.Range("C18").Value = 0.000007407357063418041
$ws.Range("D18").Value = 0.9392484394639198
$ws.Range("B19").Value = 0.0008692897494029967
$ws.Range("C19").Value = 0.000007938281114732091
$ws.Range("D19").Value = 0.913191616510361
$ws.Range("B20").Value = 0.002204643616314547
$ws.Range("C20").Value = 0.00001994844242080893
$ws.Range("D20").Value = 0.9048375108425136
$ws.Range("B21").Value = 0.002404739439149371
$ws.Range("C21").Value = 0.00002059405452534171
$ws.Range("D21").Value = 0.8563944263594084
$ws.Range("B22").Value = 0.00077441929295573
$ws.Range("C22").Value = 0.000006003815461459456
$ws.Range("D22").Value = 0.7752667729318394
$ws.Range("B23").Value = 0.0006823294780733029
$ws.Range("C23").Value = 0.000005257106300168624
$ws.Range("D23").Value = 0.7704644851359993
$ws.Range("B24").Value = 0.00085179566144044
$ws.Range("C24").Value = 0.000006437046502585442
$ws.Range("D24").Value = 0.7557031332725964
$ws.Range("B25").Value = 0.001624199425002868
$ws.Range("C25").Value = 0.00001149743325584682
$ws.Range("D25").Value = 0.7078831009822899
$ws.Range("B26").Value = 0.0006334882725942965
$ws.Range("C26").Value = 0.000004261499562245992
$ws.Range("D26").Value = 0.6727037810493416
$ws.Range("B27").Value = 0.0007973742935408388
$ws.Range("C27").Value = 0.000005327499574244013
$ws.Range("D27").Value = 0.6681303394653714
$ws.Range("B28").Value = 0.0007002228628661856
$ws.Range("C28").Value = 0.000004490485599640246
$ws.Range("D28").Value = 0.6412937705660703
$ws.Range("B29").Value = 0.001710816382477938
$ws.Range("C29").Value = 0.000009871264714706807
$ws.Range("D29").Value = 0.5769914770402956
$ws.Range("B30").Value = 0.001695109302250371
$ws.Range("C30").Value = 0.000009574074342194346
$ws.Range("D30").Value = 0.5648057225268083
$ws.Range("B31").Value = 0.00189051426708698
$ws.Range("C31").Value = 0.00001010037550791293
$ws.Range("D31").Value = 0.5342660292892795
$ws.Range("B32").Value = 0.001885146205082942
$ws.Range("C32").Value = 0.00001006816173497853
$ws.Range("D32").Value = 0.5340785615371174
$ws.Range("B33").Value = 0.001883732160283686
$ws.Range("C33").Value = 0.00000998542233475181
$ws.Range("D33").Value = 0.5300871612898526
$ws.Range("B34").Value = 0.001626144348706933
$ws.Range("C34").Value = 0.000008219759429620549
$ws.Range("D34").Value = 0.5054753863737055
$ws.Range("B35").Value = 0.001365893343099652
$ws.Range("C35").Value = 0.000006306709162659315
$ws.Range("D35").Value = 0.4617277911573506
$ws.Range("B36").Value = 0.001367326769852579
$ws.Range("C36").Value = 0.000006258582426488509
$ws.Range("D36").Value = 0.4577239738503248
$ws.Range("B37").Value = 0.001356763477407602
$ws.Range("C37").Value = 0.000006066422000285779
$ws.Range("D37").Value = 0.4471245063197771
$ws.Range("B38").Value = 0.001525839116196714
$ws.Range("C38").Value = 0.000006179671625189902
$ws.Range("D38").Value = 0.4050015207758776
$ws.Range("B39").Value = 0.001699788328446765
$ws.Range("C39").Value = 0.000006291401458154108
$ws.Range("D39").Value = 0.3701285244088641
$ws.Range("B40").Value = 0.001281930807638139
$ws.Range("C40").Value = 0.000004600319532822764
$ws.Range("D40").Value = 0.3588586455222576
$ws.Range("B41").Value = 0.001248592013351922
$ws.Range("C41").Value = 0.000003924596956786117
$ws.Range("D41").Value = 0.3143218052669017
$ws.Range("B42").Value = 0.001444682443749558
$ws.Range("C42").Value = 0.000004527565615768173
$ws.Range("D42").Value = 0.3133952125850743
$ws.Range("B43").Value = 0.001635921780826934
$ws.Range("C43").Value = 0.000004992338956549617
$ws.Range("D43").Value = 0.3051697834859844
